$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates reflecting refreshed cryptocurrency market data.
# Each entry is Cell -> NewValue. Values that look like plain numbers
# (e.g. "1.00", "7.20") must be forced to Text so Excel keeps the
# exact display string instead of normalizing them to a number.
$updates = [ordered]@{
    'D2' = '70.982.74'
    'E2' = '  +1.88%  '
    'D3' = '3.798.87'
    'E3' = '  +0.58%  '
    'D4' = '1.00'
    'E4' = '  +0.09%  '
    'D5' = '702.22'
    'E5' = '  +5.69%  '
    'D6' = '172.96'
    'E6' = '  +4.19%  '
    'D7' = '3.799.50'
    'E7' = '  +0.65%  '
    'E9' = '  +0.18%  '
    'E10' = '  +2.07%  '
    'D11' = '7.29'
    'E11' = '  +4.39%  '
    'E12' = '  +0.70%  '
    'E13' = '  +6.82%  '
    'D14' = '36.08'
    'E14' = '  +2.45%  '
    'D15' = '4.440.66'
    'E15' = '  +0.63%  '
    'D16' = '3.811.59'
    'E16' = '  +0.74%  '
    'D17' = '70.946.87'
    'E18' = '  -0.20%  '
    'D19' = '7.20'
    'E19' = '  +1.45%  '
    'E20' = '  +0.32%  '
    'D21' = '10.98'
    'E21' = '  +11.52%  '
    'D22' = '480.79'
    'E22' = '  +2.12%  '
    'E23' = '  +0.20%  '
    'D24' = '84.14'
    'E24' = '  +2.23%  '
    'E25' = '  -0.87%  '
    'D26' = '12.30'
    'E26' = '  +0.31%  '
    'D27' = '10.54'
    'E27' = '  +2.49%  '
    'E28' = '  +1.95%  '
    'D29' = '3.950.78'
    'E29' = '  +0.66%  '
    'E30' = '  -0.08%  '
    'D31' = '3.14'
    'E31' = '  +12.98%  '
    'D32' = '7.57'
    'E32' = '  +4.11%  '
    'E33' = '  +0.24%  '
    'E34' = '  +6.54%  '
    'D35' = '29.47'
    'E35' = '  +1.91%  '
    'E36' = '  +3.15%  '
    'E37' = '  +0.07%  '
    'E38' = '  +1.83%  '
    'E39' = '  +3.98%  '
    'D40' = '6.00'
    'E40' = '  +2.65%  '
    'D41' = '2.23'
    'E41' = '  +9.81%  '
    'D42' = '0.977'
    'E42' = '  +1.81%  '
    'E43' = '  +0.09%  '
    'E44' = '  -0.02%  '
    'E45' = '  +18.01%  '
    'D46' = '164.17'
    'E46' = '  +4.69%  '
    'D47' = '48.87'
    'E47' = '  +1.89%  '
    'D48' = '44.29'
    'E48' = '  -2.57%  '
    'B49' = 'TheGraph'
    'C49' = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
    'D49' = '0.301'
    'E49' = '  +0.59%  '
    'B50' = 'ONDO'
    'C50' = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
    'D50' = '1.39'
    'E50' = '  -0.86%  '
    'B51' = 'Bittensor'
    'C51' = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
    'D51' = '414.17'
    'E51' = '  +7.40%  '
}

foreach ($cellRef in $updates.Keys) {
    $value = $updates[$cellRef]
    $range = $ws.Range($cellRef)
    $isNumericLooking = $value -match '^\s*[+-]?\d+(\.\d+)?\s*$'
    if ($isNumericLooking) {
        # Force text storage so "1.00"/"7.20" style values are not
        # coerced into numbers (which would drop the formatting).
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}
